# Rename the AHB-Diff header row from the generic "_old"/"_new" suffixes
# to the concrete format-version suffixes ("_FV2210" / "_FV2304"), turn the
# header row + data into a real Excel Table (ListObject), and freeze the
# header row in place - matching the upstream commit
# "chore: adapt column header formatting to respective input file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header names for columns A:U (row 1).
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn A1:U94 into a proper table ("Table1") now that the headers carry
# their final text, so the table's column names pick it up directly.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
